$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update image dimensions (iWidth/iHeight) which drive all the recalculated formulas
$ws.Range("B1").Value = 5
$ws.Range("B2").Value = 3.5

# Update scan resolutions (xRes/yRes)
$ws.Range("B12").Value = 150
$ws.Range("B13").Value = 150

# Update selection to match the saved view state
$ws.Range("B3").Select()
